# Update profile and banner
#
# 1) Refresh the auto "Last updated" date placeholder text carried on the
#    slide master and every slide layout (10/17/2024 -> 11/4/2024).
# 2) Tweak the wording of the "BrainHack" banner textbox on slide 1:
#      "If you are interested in "   -> "If you are interested in the "
#      " about brain,"  (x2)         -> " about the brain,"
#      "for brain."                  -> "for the brain."

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Part 1: date placeholders on the slide master + all custom layouts
# ---------------------------------------------------------------------
$sm = $p.SlideMaster
$ppPlaceholderDate = 16

for ($mi = 1; $mi -le $sm.Shapes.Count; $mi++) {
    $msh = $sm.Shapes.Item($mi)
    if ($msh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        $msh.TextFrame.TextRange.Text = "11/4/2024"
    }
}

$cls = $sm.CustomLayouts
for ($li = 1; $li -le $cls.Count; $li++) {
    $cl = $cls.Item($li)
    for ($si = 1; $si -le $cl.Shapes.Count; $si++) {
        $lsh = $cl.Shapes.Item($si)
        if ($lsh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $lsh.TextFrame.TextRange.Text = "11/4/2024"
        }
    }
}

# ---------------------------------------------------------------------
# Part 2: wording tweaks in the "BrainHack" banner textbox on slide 1
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$bannerGroup = $slide.Shapes.Item(5)       # "Group 23"
$banner = $bannerGroup.GroupItems.Item(2)  # "TextBox 14"
$tr = $banner.TextFrame.TextRange

# Apply edits back-to-front so earlier character offsets stay valid while
# later ones are rewritten first (text length changes otherwise shift them).
$tr.Characters(218, 10).Text = "for the brain."
$tr.Characters(121, 13).Text = " about the brain,"
$tr.Characters(103, 13).Text = " about the brain,"
$tr.Characters(1, 25).Text = "If you are interested in the "
